$wb = $excel.ActiveWorkbook

# --- 1) Rename "Sheet2" -> "메뉴테이블" and update its selection -------------
$menuSheet = $wb.Worksheets.Item(2)
$menuSheet.Name = "메뉴테이블"
$menuSheet.Select() | Out-Null
$menuSheet.Range("C4:L13").Select() | Out-Null

# --- 2) Create the new "댓글테이블" sheet by copying "테이블 명세서" --------
#        (same layout/formatting template used by the other spec sheets)
#        and place it right after "메뉴테이블".
$specSheet = $wb.Worksheets.Item(1)
$specSheet.Copy([System.Reflection.Missing]::Value, $menuSheet) | Out-Null

$commentSheet = $wb.Worksheets.Item(3)
$commentSheet.Name = "댓글테이블"

# --- 3) Replace the table-spec body (rows 4-9) with the comment-table data -
$commentSheet.Range("C4:L13").ClearContents() | Out-Null

# Row 4 - ID (primary key)
$commentSheet.Cells.Item(4,2).Value = "댓글테이블"
$commentSheet.Cells.Item(4,3).Value = "ID"
$commentSheet.Cells.Item(4,4).Value = "BIGINT"
$commentSheet.Cells.Item(4,5).Value = "CMT_ID"
$commentSheet.Cells.Item(4,6).Value = "BIGINT"
$commentSheet.Cells.Item(4,7).Value = "NOT NULL"
$commentSheet.Cells.Item(4,8).Value = "PRIMARY KEY"

# Row 5 - 게시글ID (post id, foreign key like column)
$commentSheet.Cells.Item(5,2).Value = "tbl_comment"
$commentSheet.Cells.Item(5,3).Value = "게시글ID"
$commentSheet.Cells.Item(5,4).Value = "BIGINT"
$commentSheet.Cells.Item(5,5).Value = "CMT_P_ID"
$commentSheet.Cells.Item(5,6).Value = "BIGINT"
$commentSheet.Cells.Item(5,7).Value = "NOT NULL"

# Row 6 - 작성자
$commentSheet.Cells.Item(6,3).Value = "작성자"
$commentSheet.Cells.Item(6,4).Value = "문자열(20)"
$commentSheet.Cells.Item(6,5).Value = "CMT_WRITER"
$commentSheet.Cells.Item(6,6).Value = "VARCHAR(20)"
$commentSheet.Cells.Item(6,7).Value = "NOT NULL"

# Row 7 - 날짜
$commentSheet.Cells.Item(7,3).Value = "날짜"
$commentSheet.Cells.Item(7,4).Value = "문자열(10)"
$commentSheet.Cells.Item(7,5).Value = "CMT_DATE"
$commentSheet.Cells.Item(7,6).Value = "VARCHAR(10)"

# Row 8 - 시간
$commentSheet.Cells.Item(8,3).Value = "시간"
$commentSheet.Cells.Item(8,4).Value = "문자열(10)"
$commentSheet.Cells.Item(8,5).Value = "CMT_TIME"
$commentSheet.Cells.Item(8,6).Value = "VARCHAR(10)"

# Row 9 - 댓글
$commentSheet.Cells.Item(9,3).Value = "댓글"
$commentSheet.Cells.Item(9,4).Value = "문자열(400)"
$commentSheet.Cells.Item(9,5).Value = "CMT_TEXT"
$commentSheet.Cells.Item(9,6).Value = "VARCHAR(400)"

# --- 4) Make the new sheet the active / selected tab ------------------------
$commentSheet.Select() | Out-Null
$commentSheet.Range("M9").Select() | Out-Null
